$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Volume by Month + Device" ---
# Rename header labels (shared strings are de-duplicated automatically,
# so these renames also take care of every other worksheet that reuses
# the same labels).
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Value = "device_category"   # was "dim_deviceCategory"
$ws1.Range("D1").Value = "qty"                # was "QTY"
$ws1.Range("G1").Value = "ecr"                # was "ECR"

# --- Sheet 2: "Month Over Month Comparison" ---
$ws2 = $wb.Worksheets.Item(2)
# Add a header label in A1 (previously blank) matching the style used by
# the other header cells on row 1.
$ws2.Range("B1").Copy($ws2.Range("A1"))
$ws2.Range("A1").Value = "metric"
# The metric-name labels in column A reuse the same shared strings as the
# other sheets, so keep them in sync with the renamed labels.
$ws2.Range("A3").Value = "qty"    # was "QTY"
$ws2.Range("A6").Value = "ecr"    # was "ECR"
# Rename the row label that used to read "addsToCart".
$ws2.Range("A7").Value = "adds_to_cart"

# --- Sheet 3: "Ave Volume By Weekday" ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C1").Value = "qty"   # was "QTY" -> reuses the "qty" shared string
$ws3.Range("F1").Value = "ecr"   # was "ECR" -> reuses the "ecr" shared string

# --- Sheet 4: "Volume by Browser" ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A1").Value = "browser"   # was "dim_browser"
$ws4.Range("E1").Value = "ecr"       # was "ECR"

# Reorder the B/C/D columns (all_sessions, transactions, qty) ->
# (transactions, qty, all_sessions) for the header and every data row.
$ws4.Range("B1").Value = "transactions"
$ws4.Range("C1").Value = "qty"
$ws4.Range("D1").Value = "all_sessions"

$dataRng = $ws4.Range("B2:D58")
$vals = $dataRng.Value2
$rowCount = $vals.GetLength(0)
$newVals = New-Object 'object[,]' $rowCount,3
for ($i = 1; $i -le $rowCount; $i++) {
    $oldB = $vals[$i,1]
    $oldC = $vals[$i,2]
    $oldD = $vals[$i,3]
    $newVals[$i-1,0] = $oldC
    $newVals[$i-1,1] = $oldD
    $newVals[$i-1,2] = $oldB
}
$dataRng.Value2 = $newVals
